# "Fired logic and interface" - remove the 5 completed/dropped to-do items
# from the "Функционалности за довършване" (Features to finish) block on
# the first sheet (rows 13-17), and reflect the resulting UI state
# (active cell / minimized window) in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the five obsolete rows (old rows 13-17). Everything below shifts
# up automatically, row numbers / shared references / the sheet
# dimension are recalculated by Excel on save.
$ws.Range("A13:F17").EntireRow.Delete()

# Put the active sheet / selection where the author left it after the
# edit (first row of the now-shortened list).
$ws.Activate()
$ws.Range("B13").Select()

# The workbook window was left minimized.
$win = $wb.Windows.Item(1)
$win.WindowState = [Microsoft.Office.Interop.Excel.XlWindowState]::xlMinimized
